$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3-5 down to 4-6
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new match data
$ws.Range("A3").Value = "W2Rn64T7"
$ws.Range("B3").Value = "18/11/2024"
$ws.Range("C3").Value = "21:30"
$ws.Range("D3").Value = "ARGENTINA - TORNEO BETANO"
$ws.Range("E3").Value = "Instituto"
$ws.Range("F3").Value = "Argentinos Jrs"
$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 9.5
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 8
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 41
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 201
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 51
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 301
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 5
$ws.Range("AX3").Value = 21
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 401
$ws.Range("BC3").Value = 126
$ws.Range("BD3").Value = 126

# Row 2 (Atl. Tucuman vs Huracan) odds refresh
$ws.Range("G2").Value = 2.25
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 3.1
$ws.Range("AA2").Value = 23
$ws.Range("AE2").Value = 19
$ws.Range("AG2").Value = 8
$ws.Range("AI2").Value = 13

# Row 4 (was row 3 - Botafogo SP vs Avai) odds refresh
$ws.Range("G4").Value = 2.7
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 2.88
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 2.1
$ws.Range("Y4").Value = 12
$ws.Range("Z4").Value = 29
$ws.Range("AA4").Value = 29
$ws.Range("AC4").Value = 5.5
$ws.Range("AH4").Value = 12
$ws.Range("AJ4").Value = 29
$ws.Range("AN4").Value = 4.5

# Row 6 (was row 5 - Cerro Largo vs Wanderers) odds refresh
$ws.Range("G6").Value = 1.91
$ws.Range("I6").Value = 4.5
$ws.Range("O6").Value = 1.4
$ws.Range("P6").Value = 2.75
$ws.Range("X6").Value = 8
$ws.Range("Z6").Value = 15
$ws.Range("AG6").Value = 11
$ws.Range("AJ6").Value = 51

Write-Host "Applied weekly odds update"
